$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Config")

# Insert a new row at position 11 (pushes DeviceDiscoveryMap etc. down by one)
$ws.Rows.Item(11).Insert()

# Row 11 (new): DeviceDiscoveryIncludeAPs
$ws.Range("A11").Value = "DeviceDiscoveryIncludeAPs"
$ws.Range("B11").Value = $true
$ws.Range("C11").Value = "Add APs from CDP/LLDP discovery (SSH Only)"

# Row 10: DeviceDiscoveryIncludePhones - description text changed
$ws.Range("C10").Value = "Add phones from CDP/LLDP discovery (SNMP/SSH)"

# Restore the sheet's active selection/cursor position
[void]$ws.Range("C19").Select()
